$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Content")

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "D365"
}

$ws.Activate() | Out-Null
$ws.Range("E29").Select() | Out-Null
